$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bulk-updated
# from 45182 (2023-09-13) to 45184 (2023-09-15) for every data row (2-423).
$ws.Range("C2:C423").Value = 45184
